$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $savedStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $savedStyle
}

# Row 2
Set-TextValue "D2" '25.765.27'
Set-TextValue "E2" '  +3.22%  '

# Row 3
Set-TextValue "D3" '1.678.29'
Set-TextValue "E3" '  +2.28%  '

# Row 4
Set-TextValue "D4" '0.9980'
Set-TextValue "E4" '  -0.08%  '

# Row 5
Set-TextValue "D5" '237.32'
Set-TextValue "E5" '  +1.86%  '

# Row 6
Set-TextValue "D6" '0.9991'
Set-TextValue "E6" '  -0.09%  '

# Row 7
Set-TextValue "D7" '0.4621'
Set-TextValue "E7" '  -2.85%  '

# Row 8
Set-TextValue "E8" '  +0.23%  '

# Row 9
Set-TextValue "E9" '  +1.02%  '

# Row 10
Set-TextValue "D10" '1.672.34'
Set-TextValue "E10" '  +1.83%  '

# Row 11
Set-TextValue "D11" '0.06997'
Set-TextValue "E11" '  -0.59%  '

# Row 12
Set-TextValue "D12" '14.91'
Set-TextValue "E12" '  +2.42%  '

# Row 13
Set-TextValue "B13" 'Polkadot'
Set-TextValue "C13" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D13" '4.353'
Set-TextValue "E13" '  +0.65%  '

# Row 14
Set-TextValue "B14" 'Polygon'
Set-TextValue "C14" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D14" '0.5795'
Set-TextValue "E14" '  -1.65%  '

# Row 15
Set-TextValue "D15" '75.38'
Set-TextValue "E15" '  +2.32%  '

# Row 16
Set-TextValue "D16" '0.9989'
Set-TextValue "E16" '  -0.11%  '

# Row 17
Set-TextValue "D17" '0.9992'
Set-TextValue "E17" '  -0.05%  '

# Row 18
Set-TextValue "D18" '25.774.70'
Set-TextValue "E18" '  +3.30%  '

# Row 19
Set-TextValue "D19" '0.000006700'
Set-TextValue "E19" '  +1.97%  '

# Row 20
Set-TextValue "D20" '11.43'

# Row 21
Set-TextValue "D21" '1.884.70'
Set-TextValue "E21" '  +1.35%  '

# Row 22
Set-TextValue "D22" '4.471'
Set-TextValue "E22" '  +3.76%  '

# Row 23
Set-TextValue "D23" '8.664'
Set-TextValue "E23" '  +1.46%  '

# Row 24
Set-TextValue "D24" '5.240'
Set-TextValue "E24" '  +0.10%  '

# Row 25
Set-TextValue "D25" '134.22'
Set-TextValue "E25" '  +0.34%  '

# Row 26
Set-TextValue "D26" '15.00'

# Row 27
Set-TextValue "D27" '1.384'
Set-TextValue "E27" '  +0.36%  '

# Row 28
Set-TextValue "D28" '1.726'
Set-TextValue "E28" '  +5.27%  '

# Row 29
Set-TextValue "D29" '104.71'
Set-TextValue "E29" '  +0.44%  '

# Row 30
Set-TextValue "E30" '  +1.61%  '

# Row 31
Set-TextValue "B31" 'Filecoin'
Set-TextValue "C31" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D31" '3.618'
Set-TextValue "E31" '  +1.27%  '

# Row 32
Set-TextValue "B32" 'Stellar'
Set-TextValue "C32" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D32" '0.07681'
Set-TextValue "E32" '  +1.15%  '

# Row 33
Set-TextValue "D33" '0.04349'
Set-TextValue "E33" '  +1.88%  '

# Row 34
Set-TextValue "E34" '  +1.04%  '

# Row 35
Set-TextValue "D35" '0.6124'
Set-TextValue "E35" '  +3.22%  '

# Row 36
Set-TextValue "D36" '0.9529'
Set-TextValue "E36" '  +2.67%  '

# Row 37
Set-TextValue "D37" '0.9348'
Set-TextValue "E37" '  +8.45%  '

# Row 38
Set-TextValue "B38" 'Quant'
Set-TextValue "C38" 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D38" '108.74'
Set-TextValue "E38" '  +10.21%  '

# Row 39
Set-TextValue "B39" 'MXToken'
Set-TextValue "C39" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D39" '2.435'
Set-TextValue "E39" '  -5.24%  '

# Row 40
Set-TextValue "D40" '0.9972'
Set-TextValue "E40" '  -0.26%  '

# Row 41
Set-TextValue "D41" '1.871'
Set-TextValue "E41" '  +5.95%  '

# Row 42
Set-TextValue "D42" '0.01453'
Set-TextValue "E42" '  -2.83%  '

# Row 43
Set-TextValue "D43" '5.073'
Set-TextValue "E43" '  +8.71%  '

# Row 44
Set-TextValue "D44" '0.3730'
Set-TextValue "E44" '  +0.78%  '

# Row 45
Set-TextValue "D45" '0.1118'
Set-TextValue "E45" '  +1.59%  '

# Row 46
Set-TextValue "D46" '0.05294'
Set-TextValue "E46" '  +1.79%  '

# Row 47
Set-TextValue "B47" 'Elrond'
Set-TextValue "C47" 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue "D47" '31.34'
Set-TextValue "E47" '  +9.07%  '

# Row 48
Set-TextValue "B48" 'Aptos'
Set-TextValue "C48" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D48" '6.155'
Set-TextValue "E48" '  +0.73%  '

# Row 49
Set-TextValue "D49" '7.630'
Set-TextValue "E49" '  +6.54%  '

# Row 50
Set-TextValue "D50" '1.211'
Set-TextValue "E50" '  +3.22%  '

# Row 51
Set-TextValue "D51" '0.9996'
Set-TextValue "E51" '  -0.16%  '
